$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the pre-existing red-font style (currently sitting on empty C4) ---
# and stamp it onto C7/C8 ("NA" cells) before C4's own formatting is cleared.
$ws.Range("C4").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C4").ClearFormats()

# --- Row 2: filename ref changed from "7au_eos_5A" to "7au_eos_5A-1" ---
$ws.Range("E2").Value = "7au_eos_5A-1"

# --- Row 3: clear the Original protein / Substrate labels (now blank) ---
$ws.Range("A3:B3").ClearContents()

# --- Row 4: was a stray empty/red-font cell; now becomes a full data row ---
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 14
$ws.Range("E4").Value = "7au_eos_8A-1"
$ws.Range("F4").Value = "residues within 8A"

# --- Row 5 (new) ---
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = "7au_eos_8A-2"
$ws.Range("F5").Value = "residues within 8A"

# --- Row 6 (new) ---
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 14
$ws.Range("E6").Value = "7au_eos_8A-3"
$ws.Range("F6").Value = "residues within 8A"

# --- Row 7 (new) ---
$ws.Range("C7").Value = "NA"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = "7au_eos_mpnn-1"
$ws.Range("F7").Value = "protein mpnn entire"

# --- Row 8 (new) ---
$ws.Range("C8").Value = "NA"
$ws.Range("D8").Value = 12
$ws.Range("E8").Value = "7au_eos_mpnn-2"
$ws.Range("F8").Value = "protein mpnn entire"

# --- Dates for the new rows, copied from the existing date-formatted cell ---
$ws.Range("G2").Copy()
$ws.Range("G4:G8").PasteSpecial(-4122)
$ws.Range("G4").Value = 45446
$ws.Range("G5").Value = 45446
$ws.Range("G6").Value = 45446
$ws.Range("G7").Value = 45446
$ws.Range("G8").Value = 45446

# --- Highlight fill for A3:B8 (new light accent fill) ---
$ws.Range("A3:B8").Interior.Color = 13431551

# --- Selection, to match final saved state ---
$ws.Range("C8").Select() | Out-Null
